$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 1837
    5  = 799
    7  = 108
    13 = 135
    14 = 152
    15 = 4312
    18 = 468
    21 = 1626
    22 = 363
    23 = 44
    24 = 6
    26 = 2013
    29 = 3
    31 = 61
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
